$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows to write: id, user_id, txn_type, category, amount(text), txn_date(serial), note(text)
$rows = @(
    @(13, 2, "Expense", "bakary",      "500.00",    45962, "for morning tea"),
    @(14, 2, "Expense", "sugar",       "200.00",    45963, ""),
    @(15, 2, "Income",  "gfhg",        "54353.00",  45967, ""),
    @(16, 2, "Expense", "play",        "5000.00",   45962, ""),
    @(17, 2, "Expense", "glay",        "1230.00",   45989, ""),
    @(32, 2, "Expense", "food",        "200.00",    45988, "it was nessesary to buy"),
    @(33, 2, "Expense", "Electricity", "12000.00",  38936, ""),
    @(34, 2, "Income",  "food",        "123000.00", 45988, "")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Amount column must stay textual (e.g. "500.00" keeps trailing zeros) -
    # force text interpretation via the quote-prefix trick.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]

    # Date column: numeric serial date, formatted as a date.
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD"

    # Note column: textual, may be an explicit empty string.
    $note = $row[6]
    if ($note -eq "") {
        $ws.Cells.Item($r, 7).Value = "'"
    } else {
        $ws.Cells.Item($r, 7).NumberFormat = "@"
        $ws.Cells.Item($r, 7).Value = $note
    }

    $r++
}
